$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for week 50 (column BA), matching style/text-type of existing week-number headers
$ws.Range("BA1").Value = "'50"
$ws.Range("BA1").Font.Bold = $true
$ws.Range("BA1").HorizontalAlignment = -4108  # xlCenter

# Backfill week 49 (column AZ) values for rows that previously had no reported data for that week
$ws.Range("AZ30").Value = 1
$ws.Range("AZ53").Value = 0

# New BA (week 50) values
$ws.Range("BA2").Value = 0
$ws.Range("BA3").Value = 0
$ws.Range("BA5").Value = 0
$ws.Range("BA6").Value = 8
$ws.Range("BA7").Value = 0
$ws.Range("BA8").Value = 0
$ws.Range("BA12").Value = 0
$ws.Range("BA13").Value = 0
$ws.Range("BA14").Value = 0
$ws.Range("BA15").Value = 0
$ws.Range("BA16").Value = 0
$ws.Range("BA17").Value = 0
$ws.Range("BA19").Value = 0
$ws.Range("BA23").Value = 0
$ws.Range("BA25").Value = 0
$ws.Range("BA28").Value = 0
$ws.Range("BA29").Value = 0
$ws.Range("BA30").Value = 0
$ws.Range("BA31").Value = 0
$ws.Range("BA35").Value = 0
$ws.Range("BA36").Value = 1
$ws.Range("BA38").Value = 0
$ws.Range("BA41").Value = 0
$ws.Range("BA42").Value = 0
$ws.Range("BA43").Value = 0
$ws.Range("BA46").Value = 0
$ws.Range("BA47").Value = 0
$ws.Range("BA48").Value = 0
$ws.Range("BA49").Value = 0
$ws.Range("BA50").Value = 0
$ws.Range("BA51").Value = 0
$ws.Range("BA53").Value = 0
$ws.Range("BA54").Value = 0
$ws.Range("BA55").Value = 0
$ws.Range("BA56").Value = 0
$ws.Range("BA57").Value = 0
$ws.Range("BA58").Value = 0
$ws.Range("BA59").Value = 0
